$wb = $excel.ActiveWorkbook
$wsE = $wb.Worksheets.Item("Euramet")
$wsI = $wb.Worksheets.Item("Istruzioni Uso")

# --- Euramet sheet: calibration table values (Q3/Q1 measurement rows) ---
$wsE.Range("F7").Value = -156.5579745837849
$wsE.Range("G7").Value = -1.1772
$wsE.Range("F8").Value = -156.5579745837849
$wsE.Range("F10").Value = -157.0411782090434
$wsE.Range("G10").Value = -1.1772
$wsE.Range("G11").Value = -1.1772
$wsE.Range("E12").Value = 0
$wsE.Range("G12").Value = -1.1772
$wsE.Range("E13").Value = 0
$wsE.Range("E14").Value = 0
$wsE.Range("G14").Value = -1.1772
$wsE.Range("E15").Value = 0
$wsE.Range("F15").Value = -156.5579745837849
$wsE.Range("G15").Value = -1.1772
$wsE.Range("E16").Value = 0
$wsE.Range("F16").Value = -157.0411782090434
$wsE.Range("G16").Value = -1.1772
$wsE.Range("E17").Value = 0
$wsE.Range("G17").Value = -1.1772
$wsE.Range("E18").Value = 0
$wsE.Range("G18").Value = -1.1772
$wsE.Range("E19").Value = 0
$wsE.Range("G19").Value = -1.1772
$wsE.Range("D20").Value = 837
$wsE.Range("E20").Value = 0
$wsE.Range("F20").Value = -157.0411782090434
$wsE.Range("G20").Value = -1.3734
$wsE.Range("H20").Value = 1
$wsE.Range("D21").Value = 837
$wsE.Range("E21").Value = 0
$wsE.Range("F21").Value = -157.0411782090434
$wsE.Range("G21").Value = -1.3734
$wsE.Range("H21").Value = 1
$wsE.Range("D22").Value = 837
$wsE.Range("E22").Value = 0
$wsE.Range("F22").Value = -157.0411782090434
$wsE.Range("G22").Value = -1.3734
$wsE.Range("H22").Value = 1
$wsE.Range("D23").Value = 837
$wsE.Range("E23").Value = 0
$wsE.Range("F23").Value = -156.5579745837849
$wsE.Range("G23").Value = -1.1772
$wsE.Range("H23").Value = 1
$wsE.Range("D24").Value = 837
$wsE.Range("E24").Value = 0
$wsE.Range("F24").Value = -157.0411782090434
$wsE.Range("G24").Value = -1.3734
$wsE.Range("H24").Value = 1
$wsE.Range("D25").Value = 837
$wsE.Range("E25").Value = 0
$wsE.Range("F25").Value = -157.0411782090434
$wsE.Range("G25").Value = -1.1772
$wsE.Range("H25").Value = 1
$wsE.Range("F29").Value = -156.5579745837849
$wsE.Range("G30").Value = -1.1772
$wsE.Range("G31").Value = -1.1772
$wsE.Range("F32").Value = -157.0411782090434
$wsE.Range("G33").Value = -1.1772
$wsE.Range("G35").Value = -1.1772
$wsE.Range("G36").Value = -1.1772
$wsE.Range("F37").Value = -157.0411782090434
$wsE.Range("G37").Value = -1.1772
$wsE.Range("G38").Value = -1.1772
$wsE.Range("G39").Value = -1.1772
$wsE.Range("F40").Value = -157.0411782090434
$wsE.Range("D42").Value = 837
$wsE.Range("E42").Value = 0
$wsE.Range("F42").Value = -156.5579745837849
$wsE.Range("G42").Value = -1.1772
$wsE.Range("H42").Value = 1
$wsE.Range("D43").Value = 837
$wsE.Range("E43").Value = 0
$wsE.Range("F43").Value = -156.5579745837849
$wsE.Range("G43").Value = -1.1772
$wsE.Range("H43").Value = 1
$wsE.Range("D44").Value = 837
$wsE.Range("E44").Value = 0
$wsE.Range("F44").Value = -156.5579745837849
$wsE.Range("G44").Value = -1.1772
$wsE.Range("H44").Value = 1
$wsE.Range("D45").Value = 837
$wsE.Range("E45").Value = 0
$wsE.Range("F45").Value = -157.0411782090434
$wsE.Range("G45").Value = -1.1772
$wsE.Range("H45").Value = 1
$wsE.Range("D46").Value = 837
$wsE.Range("E46").Value = 0
$wsE.Range("F46").Value = -156.5579745837849
$wsE.Range("G46").Value = -1.3734
$wsE.Range("H46").Value = 1
$wsE.Range("D47").Value = 837
$wsE.Range("E47").Value = 0
$wsE.Range("F47").Value = -156.5579745837849
$wsE.Range("G47").Value = -1.3734
$wsE.Range("H47").Value = 1

# --- Istruzioni Uso sheet: scale / taratura parameters ---
$wsI.Range("B13").Value = 1
$wsI.Range("B20").Value = 2000

# --- Istruzioni Uso sheet: job/report header info reset to placeholders ---
$wsI.Range("B63").Value = "-"
$wsI.Range("B64").Value = "-"
$wsI.Range("B65").Value = "-"
$wsI.Range("B66").Value = "-"
$wsI.Range("B67").Value = "-"
$wsI.Range("B68").Value = "-"
